$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in existing header
$ws.Range("C2").Value = "PassingYards"

# New header cells for the three added stat columns
$ws.Range("M2").Value = "Sacks"
$ws.Range("N2").Value = "Tackles"
$ws.Range("O2").Value = "Interceptions"

# Match formatting of the other header cells (bold Times New Roman 8,
# centered horizontally+vertically, but without wrap) by copying the
# format of an existing centered header and then disabling wrap.
$ws.Range("A2").Copy() | Out-Null
$headerRange = $ws.Range("M2:O2")
$headerRange.PasteSpecial(-4122)
$headerRange.WrapText = $false

# Fill the new columns with 0 for every data row (rows 3-19).
# Clear first so the (pre-existing, empty) M19 placeholder cell drops
# back to the column's default style instead of keeping its old wrap
# formatting.
$dataRange = $ws.Range("M3:O19")
$dataRange.Clear()
for ($r = 3; $r -le 19; $r++) {
    $ws.Cells.Item($r, 13).Value = 0
    $ws.Cells.Item($r, 14).Value = 0
    $ws.Cells.Item($r, 15).Value = 0
}

# Match the plain (non-wrapped, non-bold) formatting used elsewhere
$dataRange.Font.Name = "Times New Roman"
$dataRange.Font.Size = 8

# Restore selection to match the saved workbook state
[void]$ws.Range("K21").Select()
